$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '42.049.69'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.45%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.303.83'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '315.97'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.98%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '104.42'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.60%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.622'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.15%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.26%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '39.71'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.87%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.73%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.46'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.64%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.10%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.69%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.45'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.89%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.651.60'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.80%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.303.63'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.13%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '42.029.14'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.45%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.71'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.39%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.18%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '288.29'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +11.91%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '73.92'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.12%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.32%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.30'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.00'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +6.74%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.50%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 2)
$cell.NumberFormat = "@"
$cell.Value = 'LEO'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.99'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.38%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Cosmos'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.94'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.49%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 2)
$cell.NumberFormat = "@"
$cell.Value = 'EthereumClassic'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.65'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.10%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Toncoin'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.23'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.43%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '165.59'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.19%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 2)
$cell.NumberFormat = "@"
$cell.Value = 'InjectiveProtocol'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '35.57'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.00%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Hedera'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0883'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.25%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 2)
$cell.NumberFormat = "@"
$cell.Value = 'WEMIXToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.92'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.00%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Filecoin'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.89'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.96%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Stellar'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.132'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.38%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Kaspa'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.118'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.07%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 2)
$cell.NumberFormat = "@"
$cell.Value = 'RenderToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.65'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.41%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 2)
$cell.NumberFormat = "@"
$cell.Value = 'LidoDAOToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.94'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +9.00%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 2)
$cell.NumberFormat = "@"
$cell.Value = 'VeChain'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0353'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.04%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = 'NEARProtocol'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.62'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.59%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = 'BitcoinSV'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '103.26'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +20.04%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 2)
$cell.NumberFormat = "@"
$cell.Value = 'ARBITRUM'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.49'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.81%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 2)
$cell.NumberFormat = "@"
$cell.Value = 'MultiversX'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '70.69'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.44%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Algorand'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.227'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.20%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 2)
$cell.NumberFormat = "@"
$cell.Value = 'FirstDigitalUSD'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Aave'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '116.99'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.16%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Celestia'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.06'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.46%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.10'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.06%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 2)
$cell.NumberFormat = "@"
$cell.Value = 'ordi'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '77.49'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +6.13%  '
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 2)
$cell.NumberFormat = "@"
$cell.Value = 'THORChain'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.34'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.52%  '
$cell.Style = "Normal"
